$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.43031725749783
$ws.Range("C2").Value = 9.625944543599644
$ws.Range("D2").Value = 13.75508891768489
$ws.Range("E2").Value = 14.45603335894962
$ws.Range("G2").Value = 38.68769144892295
$ws.Range("H2").Value = 16.79837864904862
$ws.Range("J2").Value = 8.781532708085813
$ws.Range("L2").Value = 11.84481059608656
$ws.Range("M2").Value = 17.2211721955651
$ws.Range("N2").Value = 19.35149609295941
$ws.Range("O2").Value = 26.89312523476997
$ws.Range("B3").Value = 17.08485573599801
$ws.Range("C3").Value = 9.489411025162307
$ws.Range("D3").Value = 13.77147696068225
$ws.Range("E3").Value = 14.49470588287755
$ws.Range("G3").Value = 38.73187193728938
$ws.Range("H3").Value = 16.84451031635287
$ws.Range("J3").Value = 8.784595368689303
$ws.Range("L3").Value = 11.84898693360549
$ws.Range("M3").Value = 17.15137770190502
$ws.Range("N3").Value = 19.40762887503675
$ws.Range("O3").Value = 26.95996148546493
$ws.Range("B4").Value = 16.87181367714146
$ws.Range("C4").Value = 9.403955871391167
$ws.Range("D4").Value = 13.78386064056383
$ws.Range("E4").Value = 14.51996197516296
$ws.Range("G4").Value = 38.76979946123107
$ws.Range("H4").Value = 16.87553636842249
$ws.Range("J4").Value = 8.786634868492484
$ws.Range("L4").Value = 11.85291002063595
$ws.Range("M4").Value = 17.11066522042558
$ws.Range("N4").Value = 19.44391276114461
$ws.Range("O4").Value = 27.0065236800211
$ws.Range("B5").Value = 16.78488108527464
$ws.Range("C5").Value = 9.368749089672153
$ws.Range("D5").Value = 13.78949099723197
$ws.Range("E5").Value = 14.53063468063497
$ws.Range("G5").Value = 38.78796529561879
$ws.Range("H5").Value = 16.88885864973559
$ws.Range("J5").Value = 8.787506090710174
$ws.Range("L5").Value = 11.85485130391644
$ws.Range("M5").Value = 17.0946245510272
$ws.Range("N5").Value = 19.45915705994924
$ws.Range("O5").Value = 27.02688496282124
$ws.Range("B6").Value = 16.77044244816542
$ws.Range("C6").Value = 9.362880634183023
$ws.Range("D6").Value = 13.79046118300982
$ws.Range("E6").Value = 14.53242988545739
$ws.Range("G6").Value = 38.79114520280768
$ws.Range("H6").Value = 16.89111180045638
$ws.Range("J6").Value = 8.787653182661201
$ws.Range("L6").Value = 11.85519436966675
$ws.Range("M6").Value = 17.09199455245139
$ws.Range("N6").Value = 19.46171607921167
$ws.Range("O6").Value = 27.0303496254343
$ws.Range("B7").Value = 16.87064158778237
$ws.Range("C7").Value = 9.403482578137531
$ws.Range("D7").Value = 13.78393420915499
$ws.Range("E7").Value = 14.52010436887243
$ws.Range("G7").Value = 38.77003348782416
$ws.Range("H7").Value = 16.87571328887044
$ws.Range("J7").Value = 8.786646455533354
$ws.Range("L7").Value = 11.85293481315677
$ws.Range("M7").Value = 17.11044664871249
$ws.Range("N7").Value = 19.44411649385581
$ws.Range("O7").Value = 27.00679266754035
$ws.Range("B8").Value = 17.31146380494976
$ws.Range("C8").Value = 9.57921944540934
$ws.Range("D8").Value = 13.76025786309675
$ws.Range("E8").Value = 14.4690544498579
$ws.Range("G8").Value = 38.70068013517535
$ws.Range("H8").Value = 16.81372408935948
$ws.Range("J8").Value = 8.782555784571999
$ws.Range("L8").Value = 11.8459691438642
$ws.Range("M8").Value = 17.19666960814494
$ws.Range("N8").Value = 19.37047400727203
$ws.Range("O8").Value = 26.91502240705616
$ws.Range("B9").Value = 18.16358813256808
$ws.Range("C9").Value = 9.909898698585886
$ws.Range("D9").Value = 13.73223559167794
$ws.Range("E9").Value = 14.3809037401354
$ws.Range("G9").Value = 38.6505775334454
$ws.Range("H9").Value = 16.71360692622949
$ws.Range("J9").Value = 8.775790233003345
$ws.Range("L9").Value = 11.84304923805124
$ws.Range("M9").Value = 17.38221822880443
$ws.Range("N9").Value = 19.24043542027594
$ws.Range("O9").Value = 26.77899432307539
$ws.Range("B10").Value = 18.77576692656239
$ws.Range("C10").Value = 10.14298443005685
$ws.Range("D10").Value = 13.7228470318766
$ws.Range("E10").Value = 14.32338615853562
$ws.Range("G10").Value = 38.6663224789155
$ws.Range("H10").Value = 16.65313909344765
$ws.Range("J10").Value = 8.771578067811422
$ws.Range("L10").Value = 11.84739278269169
$ws.Range("M10").Value = 17.52784078013558
$ws.Range("N10").Value = 19.15358610065249
$ws.Range("O10").Value = 26.70596250278965
$ws.Range("B11").Value = 19.05000083549895
$ws.Range("C11").Value = 10.24660672179805
$ws.Range("D11").Value = 13.72100027372758
$ws.Range("E11").Value = 14.29878441266829
$ws.Range("G11").Value = 38.68490617209718
$ws.Range("H11").Value = 16.6284756236205
$ws.Range("J11").Value = 8.769824992145328
$ws.Range("L11").Value = 11.85076388925958
$ws.Range("M11").Value = 17.59594070397812
$ws.Range("N11").Value = 19.11594775380985
$ws.Range("O11").Value = 26.67860361500412
$ws.Range("B12").Value = 19.1531384332507
$ws.Range("C12").Value = 10.28547765654989
$ws.Range("D12").Value = 13.72064862660878
$ws.Range("E12").Value = 14.28969250002782
$ws.Range("G12").Value = 38.69358367546864
$ws.Range("H12").Value = 16.61954533655253
$ws.Range("J12").Value = 8.769184471185238
$ws.Range("L12").Value = 11.85223974345115
$ws.Range("M12").Value = 17.62198040867753
$ws.Range("N12").Value = 19.10196285620341
$ws.Range("O12").Value = 26.669088113474
$ws.Range("B13").Value = 19.13095892502989
$ws.Range("C13").Value = 10.27712283878137
$ws.Range("D13").Value = 13.7207089130854
$ws.Range("E13").Value = 14.29164064209415
$ws.Range("G13").Value = 38.69164191960317
$ws.Range("H13").Value = 16.62145042814392
$ws.Range("J13").Value = 8.76932138305024
$ws.Range("L13").Value = 11.85191305056802
$ws.Range("M13").Value = 17.61636133988167
$ws.Range("N13").Value = 19.10496285296142
$ws.Range("O13").Value = 26.67109985649998
$ws.Range("B14").Value = 19.05850077263262
$ws.Range("C14").Value = 10.2498121667629
$ws.Range("D14").Value = 13.72096438124304
$ws.Range("E14").Value = 14.29803192475051
$ws.Range("G14").Value = 38.68558722006983
$ws.Range("H14").Value = 16.62773272016463
$ws.Range("J14").Value = 8.769771829185688
$ws.Range("L14").Value = 11.85088132257055
$ws.Range("M14").Value = 17.59807804849154
$ws.Range("N14").Value = 19.11479184359681
$ws.Range("O14").Value = 26.67780383179138
$ws.Range("B15").Value = 19.01402295107649
$ws.Range("C15").Value = 10.23303495356217
$ws.Range("D15").Value = 13.72116611143096
$ws.Range("E15").Value = 14.30197595528073
$ws.Range("G15").Value = 38.68209206209042
$ws.Range("H15").Value = 16.63163410864724
$ws.Range("J15").Value = 8.770050775339804
$ws.Range("L15").Value = 11.85027527163036
$ws.Range("M15").Value = 17.58691133525956
$ws.Range("N15").Value = 19.12084724825276
$ws.Range("O15").Value = 26.68202025799881
$ws.Range("B16").Value = 18.75775109384371
$ws.Range("C16").Value = 10.13616207208194
$ws.Range("D16").Value = 13.72301642772933
$ws.Range("E16").Value = 14.32502534941173
$ws.Range("G16").Value = 38.66533767339793
$ws.Range("H16").Value = 16.65480815887229
$ws.Range("J16").Value = 8.771695906586867
$ws.Range("L16").Value = 11.84720043534458
$ws.Range("M16").Value = 17.52342634676856
$ws.Range("N16").Value = 19.15608339920653
$ws.Range("O16").Value = 26.70786859888087
$ws.Range("B17").Value = 18.59937680438963
$ws.Range("C17").Value = 10.07610051723811
$ws.Range("D17").Value = 13.72477182136073
$ws.Range("E17").Value = 14.33956539769466
$ws.Range("G17").Value = 38.65798367811375
$ws.Range("H17").Value = 16.66975325379403
$ws.Range("J17").Value = 8.772746821318083
$ws.Range("L17").Value = 11.84567054047533
$ws.Range("M17").Value = 17.4849449895912
$ws.Range("N17").Value = 19.17817786372924
$ws.Range("O17").Value = 26.72522871241574
$ws.Range("B18").Value = 18.50788961518433
$ws.Range("C18").Value = 10.04132915709196
$ws.Range("D18").Value = 13.72600971359642
$ws.Range("E18").Value = 14.3480756263936
$ws.Range("G18").Value = 38.65482910450606
$ws.Range("H18").Value = 16.67861691787734
$ws.Range("J18").Value = 8.773366634557389
$ws.Range("L18").Value = 11.84492199395957
$ws.Range("M18").Value = 17.46298710540325
$ws.Range("N18").Value = 19.19106206972035
$ws.Range("O18").Value = 26.73576563026863
$ws.Range("B19").Value = 18.47684894839226
$ws.Range("C19").Value = 10.02951812414573
$ws.Range("D19").Value = 13.72646806756415
$ws.Range("E19").Value = 14.35098233883516
$ws.Range("G19").Value = 38.65394573247898
$ws.Range("H19").Value = 16.68166395950081
$ws.Range("J19").Value = 8.77357913346056
$ws.Range("L19").Value = 11.84469115804741
$ws.Range("M19").Value = 17.45558315160757
$ws.Range("N19").Value = 19.19545470738565
$ws.Range("O19").Value = 26.73942797773737
$ws.Range("B20").Value = 18.61627757302591
$ws.Range("C20").Value = 10.08251767822602
$ws.Range("D20").Value = 13.72456134120976
$ws.Range("E20").Value = 14.33800235653538
$ws.Range("G20").Value = 38.65865524910448
$ws.Range("H20").Value = 16.66813462067265
$ws.Range("J20").Value = 8.772633361327763
$ws.Range("L20").Value = 11.84581980956027
$ws.Range("M20").Value = 17.48902333249409
$ws.Range("N20").Value = 19.17580765625751
$ws.Range("O20").Value = 26.72332357342468
$ws.Range("B21").Value = 19.07980347554717
$ws.Range("C21").Value = 10.25784414726369
$ws.Range("D21").Value = 13.7208799161355
$ws.Range("E21").Value = 14.2961485689012
$ws.Range("G21").Value = 38.68732114069272
$ws.Range("H21").Value = 16.62587634984677
$ws.Range("J21").Value = 8.769638889976861
$ws.Range("L21").Value = 11.85117896850505
$ws.Range("M21").Value = 17.60344158407528
$ws.Range("N21").Value = 19.11189756682018
$ws.Range("O21").Value = 26.67581177552943
$ws.Range("B22").Value = 19.37857172397374
$ws.Range("C22").Value = 10.37027224694675
$ws.Range("D22").Value = 13.72049990518639
$ws.Range("E22").Value = 14.27010145435497
$ws.Range("G22").Value = 38.71561452849729
$ws.Range("H22").Value = 16.60064357804893
$ws.Range("J22").Value = 8.767817765414817
$ws.Range("L22").Value = 11.8558424222935
$ws.Range("M22").Value = 17.67968063586416
$ws.Range("N22").Value = 19.07168988627236
$ws.Range("O22").Value = 26.64968423831109
$ws.Range("B23").Value = 19.21952604535407
$ws.Range("C23").Value = 10.31047173483442
$ws.Range("D23").Value = 13.72051768468714
$ws.Range("E23").Value = 14.28388390564089
$ws.Range("G23").Value = 38.69964032214261
$ws.Range("H23").Value = 16.61389241749652
$ws.Range("J23").Value = 8.768777334259187
$ws.Range("L23").Value = 11.85324768661736
$ws.Range("M23").Value = 17.638861888502
$ws.Range("N23").Value = 19.0930069371162
$ws.Range("O23").Value = 26.66317799840384
$ws.Range("B24").Value = 18.60863808511558
$ws.Range("C24").Value = 10.07961722927716
$ws.Range("D24").Value = 13.72465578687283
$ws.Range("E24").Value = 14.33870853717298
$ws.Range("G24").Value = 38.65834828833533
$ws.Range("H24").Value = 16.6688655589184
$ws.Range("J24").Value = 8.772684607899784
$ws.Range("L24").Value = 11.8457519168239
$ws.Range("M24").Value = 17.48717899681149
$ws.Range("N24").Value = 19.17687866066305
$ws.Range("O24").Value = 26.72418315384617
$ws.Range("B25").Value = 17.93506089655152
$ws.Range("C25").Value = 9.822077835351577
$ws.Range("D25").Value = 13.73784746164813
$ws.Range("E25").Value = 14.40347515194868
$ws.Range("G25").Value = 38.65491001635363
$ws.Range("H25").Value = 16.73839386818306
$ws.Range("J25").Value = 8.777486766856068
$ws.Range("L25").Value = 11.84269512026466
$ws.Range("M25").Value = 17.33033383953335
$ws.Range("N25").Value = 19.27408305388925
$ws.Range("O25").Value = 26.81107635660421
